# Insert a new weekly data row at row 102 (pushing existing rows 102:112 down
# to 103:113), and populate it with the new Macroferia Regional de Talca -
# Esparragos record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 102; this shifts rows 102-112 down to
# 103-113 and copies formatting (incl. the date number format) from the row
# above, matching the existing sheet's style usage.
$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value2 = 5
$ws.Range("B102").Value2 = "Macroferia Regional de Talca"
$ws.Range("C102").Value2 = "Maule"
$ws.Range("D102").Value2 = 45209
$ws.Range("E102").Value2 = 7
$ws.Range("F102").Value2 = 300000000
$ws.Range("G102").Value2 = "Espárragos"
$ws.Range("H102").Value2 = "Verde"
$ws.Range("I102").Value2 = "Primera"
$ws.Range("J102").Value2 = 5000
$ws.Range("K102").Value2 = 1100
$ws.Range("L102").Value2 = 1200
$ws.Range("M102").Value2 = 1140
$ws.Range("N102").Value2 = "`$/kilo"
$ws.Range("O102").Value2 = "Provincia de Linares"
$ws.Range("P102").Value2 = 1140
$ws.Range("Q102").Value2 = 1
$ws.Range("R102").Value2 = "Hortaliza"
